$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("C2").Value = 71
$ws.Range("G2").Value = 87
$ws.Range("H2").Value = 109
$ws.Range("B3").Value = 79
$ws.Range("F3").Value = 142
$ws.Range("I3").Value = 197
$ws.Range("B6").Value = 385
$ws.Range("D6").Value = 427
$ws.Range("F6").Value = 559
$ws.Range("H6").Value = 456
$ws.Range("K6").Value = 519
$ws.Range("B7").Value = 521
$ws.Range("C7").Value = 649
$ws.Range("D7").Value = 667
$ws.Range("F7").Value = 808
$ws.Range("G7").Value = 675
$ws.Range("H7").Value = 741
$ws.Range("I7").Value = 851
$ws.Range("K7").Value = 918

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("B3").Value = 1
$ws.Range("B7").Value = 17

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("B5").Value = 6
$ws.Range("B6").Value = 8

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("F10").Value = 5
$ws.Range("I21").Value = 16
$ws.Range("K29").Value = 23
$ws.Range("B32").Value = 17
$ws.Range("B41").Value = 8
$ws.Range("G47").Value = 19
$ws.Range("K52").Value = 10
$ws.Range("H53").Value = 110
$ws.Range("K54").Value = 6
$ws.Range("B62").Value = 7
$ws.Range("F68").Value = 6
$ws.Range("C76").Value = 15
$ws.Range("D76").Value = 16
$ws.Range("H77").Value = 31
$ws.Range("B98").Value = 521
$ws.Range("C98").Value = 649
$ws.Range("D98").Value = 667
$ws.Range("F98").Value = 808
$ws.Range("G98").Value = 675
$ws.Range("H98").Value = 741
$ws.Range("I98").Value = 851
$ws.Range("K98").Value = 918

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("H6").Value = 71
$ws.Range("H7").Value = 110

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("C2").Value = 3
$ws.Range("D6").Value = 10
$ws.Range("C7").Value = 15
$ws.Range("D7").Value = 16

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("B3").Value = 1
$ws.Range("B7").Value = 7

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("K5").Value = 18
$ws.Range("K6").Value = 23

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("H2").Value = 6
$ws.Range("H7").Value = 31

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("F3").Value = 2
$ws.Range("F5").Value = 3
$ws.Range("F6").Value = 5

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K4").Value = 3
$ws.Range("K5").Value = 6

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K5").Value = 7
$ws.Range("K6").Value = 10

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("G2").Value = 4
$ws.Range("G6").Value = 19

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("F5").Value = 4
$ws.Range("F6").Value = 6

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("I3").Value = 6
$ws.Range("I7").Value = 16
